# Apply the "updated trends paper for new bootstrapping" edits:
# a handful of standard-error cells (shown as "(0.xxxx)" text) in columns C/D
# get new values reflecting the new bootstrapping results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Value  = "(0.0001)"
$ws.Range("C7").Value  = "(0.0003)"
$ws.Range("C9").Value  = "(0.0003)"
$ws.Range("D9").Value  = "(0.0007)"
$ws.Range("C11").Value = "(0.0)"
$ws.Range("D11").Value = "(0.0011)"
$ws.Range("D13").Value = "(0.0007)"
$ws.Range("D15").Value = "(0.0009)"
